$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Helper: write a "DD-MM-YYYY" looking string into a cell as literal TEXT.
# Excel auto-detects ambiguous day<=12 strings like "01-11-2021" as dates and
# silently rewrites them into a date serial number (with a new NumberFormat
# style) when assigned straight to .Value / .Value2 / .Formula. To avoid that,
# build the text via a formula result (so it is never parsed as a typed
# literal), copy it, and paste-special VALUES ONLY into the destination cell.
# The scratch cell sits far outside the sheet's real data and is cleared
# immediately afterwards so it never affects the used range / dimension. ---
function Set-TextValue {
    param($cell, $text)
    $scratch = $ws.Cells.Item(2000, 26)
    $scratch.Formula = '="' + $text + '"'
    $scratch.Copy()
    $cell.PasteSpecial(-4163)
    $scratch.ClearContents()
}

# --- Row 302 (previously the last row): fill in the B/C values that were
# missing before, matching the rest of the table. ---
$ws.Cells.Item(302, 2).Value = 187
$ws.Cells.Item(302, 3).Value = 628

# --- Append the new daily rows 303-307. ---
Set-TextValue $ws.Cells.Item(303, 1) "29-10-2021"
$ws.Cells.Item(303, 2).Value = 187
$ws.Cells.Item(303, 3).Value = 628
$ws.Cells.Item(303, 4).Value = 3940
$ws.Cells.Item(303, 5).Value = 30

Set-TextValue $ws.Cells.Item(304, 1) "30-10-2021"
$ws.Cells.Item(304, 4).Value = 3940
$ws.Cells.Item(304, 5).Value = 30

Set-TextValue $ws.Cells.Item(305, 1) "31-10-2021"
$ws.Cells.Item(305, 4).Value = 3940
$ws.Cells.Item(305, 5).Value = 30

Set-TextValue $ws.Cells.Item(306, 1) "01-11-2021"
$ws.Cells.Item(306, 4).Value = 3940
$ws.Cells.Item(306, 5).Value = 30

Set-TextValue $ws.Cells.Item(307, 1) "02-11-2021"
$ws.Cells.Item(307, 4).Value = 3940
$ws.Cells.Item(307, 5).Value = 30
